# Refresh crypto price/volume figures (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage for cells whose
# content looks numeric (e.g. "28.24"), so Excel does not silently
# convert it to a Number and strip meaningful trailing zeros.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '29.715.93'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.602.08'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  +0.10%  '
Set-TextValue $ws.Range("D8") '28.24'
$ws.Range("E8").Value = '  +5.51%  '
$ws.Range("E9").Value = '  +1.81%  '
$ws.Range("E10").Value = '  +0.78%  '
Set-TextValue $ws.Range("D11") '0.0906'
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = '1.831.51'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '1.606.06'
$ws.Range("E13").Value = '  -0.11%  '
Set-TextValue $ws.Range("D14") '0.552'
$ws.Range("E14").Value = '  +2.61%  '
$ws.Range("D15").Value = '29.703.69'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("E16").Value = '  +0.62%  '
Set-TextValue $ws.Range("D17") '64.00'
Set-TextValue $ws.Range("D18") '242.43'
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("E19").Value = '  +4.75%  '
$ws.Range("D20").Value = '0.0₃0696'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("E21").Value = '  +0.24%  '
Set-TextValue $ws.Range("D22") '4.03'
$ws.Range("E22").Value = '  -0.59%  '
Set-TextValue $ws.Range("D23") '9.42'
$ws.Range("E23").Value = '  +1.76%  '
Set-TextValue $ws.Range("D24") '2.10'
$ws.Range("E24").Value = '  +0.71%  '
Set-TextValue $ws.Range("D25") '155.17'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("E29").Value = '  +0.16%  '
Set-TextValue $ws.Range("D30") '0.0477'
$ws.Range("E30").Value = '  +0.90%  '
Set-TextValue $ws.Range("D31") '1.07'
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("E32").Value = '  -0.27%  '
$ws.Range("E33").Value = '  +2.69%  '
$ws.Range("D34").Value = '1.420.55'
$ws.Range("E34").Value = '  -1.32%  '
$ws.Range("E35").Value = '  +3.70%  '
Set-TextValue $ws.Range("D36") '1.04'
$ws.Range("E36").Value = '  -0.74%  '
Set-TextValue $ws.Range("D37") '2.88'
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("E39").Value = '  +1.42%  '
$ws.Range("E40").Value = '  +1.21%  '
Set-TextValue $ws.Range("D41") '55.65'
$ws.Range("E41").Value = '  +0.75%  '
Set-TextValue $ws.Range("D42") '0.0493'
$ws.Range("E42").Value = '  +5.75%  '
Set-TextValue $ws.Range("D43") '0.817'
$ws.Range("E43").Value = '  +2.39%  '
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D46") '0.995'
$ws.Range("E46").Value = '  +18.70%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D47") '67.05'
$ws.Range("E47").Value = '  +1.54%  '
$ws.Range("E48").Value = '  +1.08%  '
$ws.Range("D49").Value = '1.741.36'
$ws.Range("E49").Value = '  -0.63%  '
Set-TextValue $ws.Range("D50") '86.66'
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("E51").Value = '  +0.68%  '
